$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column C (Förändrad) dates from 2023-10-05 (45204) to 2023-10-08 (45207)
# for rows 2-7, keeping existing number formatting/style intact.
for ($row = 2; $row -le 7; $row++) {
    $ws.Cells.Item($row, 3).Value = 45207
}
